$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.151.89"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.224.48"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "293.55"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").Value = "87.84"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "30.65"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +6.49%  "
$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "0.113"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "2.569.44"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "13.83"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "2.217.57"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "0.736"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "40.096.23"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").Value = "5.78"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "65.65"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "235.90"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "23.19"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "9.33"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  -5.73%  "
$ws.Range("D31").Value = "158.55"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  +6.91%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").Value = "0.0996"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "15.66"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "2.086.83"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "3.75"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").Value = "19.28"
$ws.Range("E44").Value = "  +9.74%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("E48").Value = "  -13.12%  "
$ws.Range("D49").Value = "2.442.10"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("E51").Value = "  +3.65%  "
